$d = $word.ActiveDocument

$old = "Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Διδύμων 2022: 14-23 Φεβρουαρίου, 14-24 Μαρτίου"
$new = "2022 Ημερομηνίες παρατήρησης για τον αστερισμό του Αστερισμός Διδύμων: 14-23 Φεβρουαρίου, 14-24 Μαρτίου"

$find = $d.Content.Find
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
